$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.002.82'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '1.716.22'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.29'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3974'
$ws.Range('E7').Value = '  +0.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4119'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.005'
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.50'
$ws.Range('E11').Value = '  +3.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08949'
$ws.Range('E12').Value = '  +1.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.744'
$ws.Range('E13').Value = '  +7.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.23'
$ws.Range('E14').Value = '  +7.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.162'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001375'
$ws.Range('E16').Value = '  +4.13%  '
$ws.Range('D17').Value = '1.694.56'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '100.65'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07161'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.22'
$ws.Range('E20').Value = '  +2.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.537'
$ws.Range('E21').Value = '  +6.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.57'
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('D24').Value = '24.992.57'
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.142'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.343'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.24'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.249'
$ws.Range('E28').Value = '  +24.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '165.34'
$ws.Range('E29').Value = '  +2.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '140.76'
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.242'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.901'
$ws.Range('E32').Value = '  +11.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09052'
$ws.Range('E33').Value = '  +5.13%  '
$ws.Range('D34').Value = '1.881.51'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.089'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.03010'
$ws.Range('E36').Value = '  +10.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2810'
$ws.Range('E37').Value = '  +2.01%  '
$ws.Range('E38').Value = '  -3.42%  '
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.62'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09318'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8124'
$ws.Range('E42').Value = '  +5.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.487'
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.72'
$ws.Range('E44').Value = '  +5.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7387'
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.657'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.274'
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.349'
$ws.Range('E49').Value = '  +1.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.38'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '93.93'
$ws.Range('E51').Value = '  +4.50%  '
